# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for each game row on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 2
    4  = 1
    5  = 2
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 2
    11 = 2
    12 = 1
    13 = 3
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 2
    22 = 0
    23 = 1
    25 = 0
    26 = 0
    27 = 2
    28 = 1
    29 = 3
    30 = 0
    31 = 0
    32 = 1
    33 = 0
    34 = 3
    35 = 2
    36 = 0
    37 = 2
    38 = 2
    39 = 0
    40 = 2
    41 = 1
    42 = 3
    43 = 1
    44 = 1
    45 = 1
    46 = 2
    47 = 2
    48 = 1
    49 = 1
    50 = 2
    51 = 1
    52 = 1
    53 = 2
    54 = 0
    55 = 2
    56 = 2
    57 = 2
    58 = 4
    59 = 0
    60 = 2
    61 = 1
    62 = 2
    63 = 1
    64 = 1
    65 = 1
    68 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
